# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N ("Late") on the
#   "Repayment schedule" sheet, shifting the "Late" / "heading" / "Outstanding"
#   columns one place to the right (N->O, O->P, P->Q).
# - Make "Repayment schedule" the active sheet / tab, with cell J18 selected.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before column 14 (N). Everything from N onward
# shifts right by one column (N->O, O->P, P->Q).
$wsSchedule.Columns.Item(14).Insert()

# The newly inserted column inherits the width of the column immediately
# to its left (column M / 13), matching Excel's default insert behaviour.
$wsSchedule.Columns.Item(14).ColumnWidth = $wsSchedule.Columns.Item(13).ColumnWidth

# Switch to the "Repayment schedule" sheet and select cell J18, leaving the
# "Transactions" sheet no longer the active tab.
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("J18").Select() | Out-Null
